$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B20").Value = "Case File – Default assignee"
$ws.Range("C20").Value = "CASE_FILE"
$ws.Range("D20").Value = "participants.?[participantType == 'assignee'].isEmpty()"
$ws.Range("G20").Value = "assignee, ann-acm"

$ws.Range("B21").Value = "Case File – Default access"
$ws.Range("C21").Value = "CASE_FILE"
$ws.Range("D21").Value = "participants.?[participantType == '*'].isEmpty()"
$ws.Range("G21").Value = "*, *"

$ws.Range("G22").Select()
